# Fix latency units in report sheet:
#  - O2 header: "Utility" -> "Utility (Percent)"
#  - Columns I, J, K (rows 3-23): append " msec" to the existing value

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header label in O2
$ws.Range("O2").Value = "Utility (Percent)"

# Append " msec" to the latency columns (I, J, K) for rows 3 through 23
for ($row = 3; $row -le 23; $row++) {
    foreach ($col in @("I", "J", "K")) {
        $cell = $ws.Range("$col$row")
        $current = $cell.Value()
        if ($current -ne $null -and -not ("$current".EndsWith(" msec"))) {
            $cell.Value = "$current msec"
        }
    }
}
